$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 25.8200000000006
$ws.Range("H2").Value = [double]"1.641734602033503e-16"
$ws.Range("K2").Value = 55.99839804036411
$ws.Range("L2").Value = "[49.90214466576804, 62.094651414960175]"
$ws.Range("O2").Value = 1.691868716347656
$ws.Range("P2").Value = "[1.5786581702723481, 1.8050792624229643]"
$ws.Range("S2").Value = 58.96487386093517
$ws.Range("T2").Value = "[55.067160395286066, 62.862587326584276]"
$ws.Range("W2").Value = 18.8674674674679
$ws.Range("X2").Value = 18.40224224224267
$ws.Range("Y2").Value = 19.33269269269314

# Row 3 updates
$ws.Range("B3").Value = 0
$ws.Range("E3").Value = 24.80000000000044
$ws.Range("H3").Value = [double]"1.641734602033503e-16"
$ws.Range("K3").Value = 53.70498598649592
$ws.Range("L3").Value = "[43.70170559125532, 63.70826638173652]"
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0.1572368695490383
$ws.Range("P3").Value = "[-0.056605273037655834, 0.37107901213573236]"
$ws.Range("Q3").Value = 0.1482721521368495
$ws.Range("R3").Value = 0.1482721521368495
$ws.Range("S3").Value = 55.56201752868665
$ws.Range("T3").Value = "[49.38899843032503, 61.73503662704826]"
$ws.Range("W3").Value = 24.17937937937981
$ws.Range("X3").Value = 23.33533533533574
$ws.Range("Y3").Value = 25.02342342342387
